$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save row 2 values (A2:E2) before overwriting
$row2 = @(
    $ws.Range("A2").Value2,
    $ws.Range("B2").Value2,
    $ws.Range("C2").Value2,
    $ws.Range("D2").Value2,
    $ws.Range("E2").Value2
)

# Save row 4 values (A4:E4) before overwriting
$row4 = @(
    $ws.Range("A4").Value2,
    $ws.Range("B4").Value2,
    $ws.Range("C4").Value2,
    $ws.Range("D4").Value2,
    $ws.Range("E4").Value2
)

# Write row 4's original data into row 2
$ws.Range("A2").Value2 = $row4[0]
$ws.Range("B2").Value2 = $row4[1]
$ws.Range("C2").Value2 = $row4[2]
$ws.Range("D2").Value2 = $row4[3]
$ws.Range("E2").Value2 = $row4[4]

# Write row 2's original data into row 4
$ws.Range("A4").Value2 = $row2[0]
$ws.Range("B4").Value2 = $row2[1]
$ws.Range("C4").Value2 = $row2[2]
$ws.Range("D4").Value2 = $row2[3]
$ws.Range("E4").Value2 = $row2[4]
